$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fecha" (date, col D), "Calidad" (col L), "Volumen" (col M), "Precio minimo"
# (col N), "Precio maximo" (col O), "Precio promedio ponderado" (col P) and
# "Precio $/Kg" (col S) values for rows 2-14 were re-shuffled across rows
# (a data correction / weekly re-sort). Apply the new value for each row.

$rows = @(
    @{ Row = 2;  D = 44447; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 3;  D = 44839; L = "Primera";  M = 120; N = 25000; O = 26000; P = 25500; S = 2550 },
    @{ Row = 4;  D = 44460; L = "Especial"; M = 60;  N = 31000; O = 32000; P = 31500; S = 3150 },
    @{ Row = 5;  D = 44460; L = "Primera";  M = 30;  N = 30000; O = 30000; P = 30000; S = 3000 },
    @{ Row = 6;  D = 44841; L = "Primera";  M = 60;  N = 23000; O = 24000; P = 23500; S = 2350 },
    @{ Row = 7;  D = 44446; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 8;  D = 44461; L = "Especial"; M = 60;  N = 31000; O = 32000; P = 31500; S = 3150 },
    @{ Row = 9;  D = 44461; L = "Primera";  M = 30;  N = 30000; O = 30000; P = 30000; S = 3000 },
    @{ Row = 10; D = 44448; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 11; D = 44848; L = "Especial"; M = 60;  N = 24000; O = 25000; P = 24500; S = 2450 },
    @{ Row = 12; D = 44848; L = "Primera";  M = 120; N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 13; D = 44452; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 14; D = 44487; L = "Primera";  M = 30;  N = 23000; O = 24000; P = 23500; S = 2350 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("D$i").Value = $r.D
    $ws.Range("L$i").Value = $r.L
    $ws.Range("M$i").Value = $r.M
    $ws.Range("N$i").Value = $r.N
    $ws.Range("O$i").Value = $r.O
    $ws.Range("P$i").Value = $r.P
    $ws.Range("S$i").Value = $r.S
}
